$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.018.76'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '1.871.59'
$ws.Range("E3").Value = '  -2.42%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'319.18"
$ws.Range("E5").Value = '  -3.25%  '
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = "'0.5049"
$ws.Range("E7").Value = '  -3.21%  '
$ws.Range("D8").Value = "'0.3961"
$ws.Range("E8").Value = '  -2.99%  '
$ws.Range("D9").Value = "'0.08209"
$ws.Range("E9").Value = '  -3.38%  '
$ws.Range("D10").Value = "'42.13"
$ws.Range("E10").Value = '  -2.39%  '
$ws.Range("E11").Value = '  -3.00%  '
$ws.Range("D12").Value = "'23.51"
$ws.Range("E12").Value = '  +4.75%  '
$ws.Range("D13").Value = '1.871.96'
$ws.Range("E13").Value = '  -2.43%  '
$ws.Range("D14").Value = "'6.291"
$ws.Range("E14").Value = '  -1.88%  '
$ws.Range("D15").Value = "'7.196"
$ws.Range("E15").Value = '  -2.74%  '
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("E18").Value = '  -2.28%  '
$ws.Range("D19").Value = "'0.06430"
$ws.Range("D20").Value = "'18.12"
$ws.Range("E20").Value = '  -0.84%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").Value = '30.006.33'
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("E23").Value = '  -2.51%  '
$ws.Range("D24").Value = "'11.13"
$ws.Range("E24").Value = '  -1.55%  '
$ws.Range("E25").Value = '  -2.73%  '
$ws.Range("D26").Value = '2.088.36'
$ws.Range("E26").Value = '  -2.48%  '
$ws.Range("D27").Value = "'21.25"
$ws.Range("E27").Value = '  +0.79%  '
$ws.Range("D28").Value = "'161.08"
$ws.Range("E28").Value = '  +0.53%  '
$ws.Range("D29").Value = "'2.220"
$ws.Range("E29").Value = '  -9.23%  '
$ws.Range("D30").Value = "'127.35"
$ws.Range("E30").Value = '  -1.33%  '
$ws.Range("D31").Value = "'1.072"
$ws.Range("E31").Value = '  -0.26%  '
$ws.Range("D32").Value = "'0.1034"
$ws.Range("D33").Value = "'5.947"
$ws.Range("D34").Value = "'3.669"
$ws.Range("E34").Value = '  +1.10%  '
$ws.Range("E35").Value = '  -2.07%  '
$ws.Range("D36").Value = "'5.221"
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("D37").Value = "'0.06356"
$ws.Range("E37").Value = '  -3.84%  '
$ws.Range("D38").Value = "'0.2143"
$ws.Range("E38").Value = '  -2.88%  '
$ws.Range("D39").Value = "'1.175"
$ws.Range("E39").Value = '  -4.41%  '
$ws.Range("D40").Value = "'8.493"
$ws.Range("E40").Value = '  -4.65%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = "'0.6309"
$ws.Range("E41").Value = '  -3.09%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = "'1.219"
$ws.Range("E42").Value = '  -2.46%  '
$ws.Range("D43").Value = "'11.30"
$ws.Range("E43").Value = '  -2.69%  '
$ws.Range("D44").Value = "'0.9999"
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = "'0.5919"
$ws.Range("E45").Value = '  -3.84%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'12.97"
$ws.Range("E46").Value = '  -1.73%  '
$ws.Range("D47").Value = "'2.097"
$ws.Range("E47").Value = '  +0.89%  '
$ws.Range("D48").Value = "'3.622"
$ws.Range("D49").Value = "'122.58"
$ws.Range("E49").Value = '  -1.40%  '
$ws.Range("E50").Value = '  -3.21%  '
$ws.Range("D51").Value = "'77.40"
$ws.Range("E51").Value = '  -3.07%  '
